$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2225
$ws.Cells.Item(40, 9).Value = 2152
$ws.Cells.Item(40, 11).Value = 2152
$ws.Cells.Item(40, 13).Value = -1977

$ws.Cells.Item(98, 8).Value = 472.36365
$ws.Cells.Item(98, 9).Value = 219.6
$ws.Cells.Item(98, 11).Value = 219.6
$ws.Cells.Item(98, 13).Value = 1278.4

$ws.Cells.Item(100, 8).Value = 1724.1305
$ws.Cells.Item(100, 9).Value = 1556.2667
$ws.Cells.Item(100, 10).Value = 2038.875
$ws.Cells.Item(100, 11).Value = 1556.2667
$ws.Cells.Item(100, 12).Value = 2038.875
$ws.Cells.Item(100, 13).Value = -1015.2667
$ws.Cells.Item(100, 14).Value = -3120.875

$ws.Cells.Item(109, 8).Value = 33000
$ws.Cells.Item(109, 10).Value = 33000
$ws.Cells.Item(109, 12).Value = 33000
$ws.Cells.Item(109, 14).Value = -35774

$ws.Cells.Item(122, 8).Value = 472.36365
$ws.Cells.Item(122, 9).Value = 219.6
$ws.Cells.Item(122, 11).Value = 658.8
$ws.Cells.Item(122, 13).Value = 1791.2

$ws.Cells.Item(138, 8).Value = 2289.8164
$ws.Cells.Item(138, 9).Value = 1854.2916
$ws.Cells.Item(138, 10).Value = 2707.92
$ws.Cells.Item(138, 11).Value = 5562.8748
$ws.Cells.Item(138, 12).Value = 8123.76
$ws.Cells.Item(138, 13).Value = -422.8747999999996
$ws.Cells.Item(138, 14).Value = -18403.76

$ws.Cells.Item(141, 8).Value = 3091.9773
$ws.Cells.Item(141, 9).Value = 879.55554
$ws.Cells.Item(141, 11).Value = 2638.66662
$ws.Cells.Item(141, 13).Value = 2541.33338

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 361678.62
$ws.Cells.Item(32, 9).Value = 386812.38
$ws.Cells.Item(32, 11).Value = 386812.38
$ws.Cells.Item(32, 13).Value = -386525.38

$ws.Cells.Item(33, 8).Value = 67521.75

$ws.Cells.Item(36, 8).Value = 36411.6
$ws.Cells.Item(36, 9).Value = 8250
$ws.Cells.Item(36, 10).Value = 55186
$ws.Cells.Item(36, 11).Value = 8250
$ws.Cells.Item(36, 12).Value = 55186
$ws.Cells.Item(36, 13).Value = -7904
$ws.Cells.Item(36, 14).Value = -55878

$ws.Cells.Item(110, 8).Value = 1629.8334
$ws.Cells.Item(110, 9).Value = 1795.8
$ws.Cells.Item(110, 10).Value = 800
$ws.Cells.Item(110, 11).Value = 1795.8
$ws.Cells.Item(110, 12).Value = 800
$ws.Cells.Item(110, 13).Value = 249.2
$ws.Cells.Item(110, 14).Value = -4890

$ws.Cells.Item(128, 8).Value = 38314.5
$ws.Cells.Item(128, 10).Value = 38314.5
$ws.Cells.Item(128, 12).Value = 38314.5
$ws.Cells.Item(128, 14).Value = -48274.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 3262.5417
$ws.Cells.Item(22, 9).Value = 3262.5417
$ws.Cells.Item(22, 11).Value = 3262.5417
$ws.Cells.Item(22, 13).Value = -3089.5417

$ws.Cells.Item(107, 8).Value = 251650
$ws.Cells.Item(107, 9).Value = 501000
$ws.Cells.Item(107, 10).Value = 2300
$ws.Cells.Item(107, 11).Value = 501000
$ws.Cells.Item(107, 12).Value = 2300
$ws.Cells.Item(107, 13).Value = -499080
$ws.Cells.Item(107, 14).Value = -6140

$ws.Cells.Item(134, 8).Value = 2127.7693
$ws.Cells.Item(134, 9).Value = 1873.88
$ws.Cells.Item(134, 11).Value = 5621.64
$ws.Cells.Item(134, 13).Value = -3086.64

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1360.4048
$ws.Cells.Item(132, 9).Value = 943.7879
$ws.Cells.Item(132, 11).Value = 2831.3637
$ws.Cells.Item(132, 13).Value = -301.3636999999999

$ws.Cells.Item(134, 8).Value = 4635.1875
$ws.Cells.Item(134, 9).Value = 4485.4443
$ws.Cells.Item(134, 10).Value = 5443.8
$ws.Cells.Item(134, 11).Value = 13456.3329
$ws.Cells.Item(134, 12).Value = 16331.4
$ws.Cells.Item(134, 13).Value = -10921.3329
$ws.Cells.Item(134, 14).Value = -21401.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 968.41174
$ws.Cells.Item(5, 9).Value = 536.8461
$ws.Cells.Item(5, 11).Value = 1610.5383
$ws.Cells.Item(5, 13).Value = -1498.5383

$ws.Cells.Item(22, 8).Value = 966.3171
$ws.Cells.Item(22, 10).Value = 960.8570999999999
$ws.Cells.Item(22, 12).Value = 2882.5713
$ws.Cells.Item(22, 14).Value = -3220.5713

$ws.Cells.Item(27, 8).Value = 966.3171
$ws.Cells.Item(27, 10).Value = 960.8570999999999
$ws.Cells.Item(27, 12).Value = 2882.5713
$ws.Cells.Item(27, 14).Value = -3086.5713

$ws.Cells.Item(32, 8).Value = 27781112
$ws.Cells.Item(32, 10).Value = 27781112
$ws.Cells.Item(32, 12).Value = 83343336
$ws.Cells.Item(32, 14).Value = -83343902

$ws.Cells.Item(39, 8).Value = 1714.1904
$ws.Cells.Item(39, 9).Value = 500.75
$ws.Cells.Item(39, 10).Value = 1999.7059
$ws.Cells.Item(39, 11).Value = 1502.25
$ws.Cells.Item(39, 12).Value = 5999.1177
$ws.Cells.Item(39, 13).Value = -1208.25
$ws.Cells.Item(39, 14).Value = -6587.1177

$ws.Cells.Item(46, 8).Value = 659.17645
$ws.Cells.Item(46, 9).Value = 389.55554
$ws.Cells.Item(46, 10).Value = 962.5
$ws.Cells.Item(46, 11).Value = 1168.66662
$ws.Cells.Item(46, 12).Value = 2887.5
$ws.Cells.Item(46, 13).Value = -1077.66662
$ws.Cells.Item(46, 14).Value = -3069.5

$ws.Cells.Item(58, 8).Value = 1933.7037
$ws.Cells.Item(58, 9).Value = 900
$ws.Cells.Item(58, 10).Value = 1973.4615
$ws.Cells.Item(58, 11).Value = 2700
$ws.Cells.Item(58, 12).Value = 5920.3845
$ws.Cells.Item(58, 13).Value = -2572
$ws.Cells.Item(58, 14).Value = -6176.3845

$ws.Cells.Item(110, 8).Value = 14674.667
$ws.Cells.Item(110, 10).Value = 14674.667
$ws.Cells.Item(110, 12).Value = 44024.001
$ws.Cells.Item(110, 14).Value = -52204.001

$ws.Cells.Item(122, 8).Value = 11124.5
$ws.Cells.Item(122, 10).Value = 18015.834
$ws.Cells.Item(122, 12).Value = 162142.506
$ws.Cells.Item(122, 14).Value = -167042.506

$ws.Cells.Item(123, 8).Value = 5426.6665
$ws.Cells.Item(123, 9).Value = 2030
$ws.Cells.Item(123, 10).Value = 7125
$ws.Cells.Item(123, 11).Value = 6090
$ws.Cells.Item(123, 12).Value = 21375
$ws.Cells.Item(123, 13).Value = -3640
$ws.Cells.Item(123, 14).Value = -26275

$ws.Cells.Item(132, 8).Value = 2318.4062
$ws.Cells.Item(132, 9).Value = 1983.3334
$ws.Cells.Item(132, 10).Value = 2449.5217
$ws.Cells.Item(132, 11).Value = 17850.0006
$ws.Cells.Item(132, 12).Value = 22045.6953
$ws.Cells.Item(132, 13).Value = -15320.0006
$ws.Cells.Item(132, 14).Value = -27105.6953

$ws.Cells.Item(135, 8).Value = 968.41174
$ws.Cells.Item(135, 9).Value = 536.8461
$ws.Cells.Item(135, 11).Value = 4831.6149
$ws.Cells.Item(135, 13).Value = -2296.6149

$ws.Cells.Item(136, 8).Value = 1828.4286
$ws.Cells.Item(136, 9).Value = 959.8
$ws.Cells.Item(136, 11).Value = 2879.4
$ws.Cells.Item(136, 13).Value = 2220.6

$ws.Cells.Item(137, 8).Value = 6181110.5
$ws.Cells.Item(137, 9).Value = 33356740
$ws.Cells.Item(137, 10).Value = 4830.727
$ws.Cells.Item(137, 11).Value = 100070220
$ws.Cells.Item(137, 12).Value = 14492.181
$ws.Cells.Item(137, 13).Value = -100065120
$ws.Cells.Item(137, 14).Value = -24692.181

$ws.Cells.Item(140, 8).Value = 1456.8125
$ws.Cells.Item(140, 9).Value = 1309.0834
$ws.Cells.Item(140, 10).Value = 1900
$ws.Cells.Item(140, 11).Value = 3927.2502
$ws.Cells.Item(140, 12).Value = 5700
$ws.Cells.Item(140, 13).Value = 1252.7498
$ws.Cells.Item(140, 14).Value = -16060

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 6950
$ws.Cells.Item(18, 10).Value = 6950
$ws.Cells.Item(18, 12).Value = 6950
$ws.Cells.Item(18, 14).Value = -7536

$ws.Cells.Item(122, 8).Value = 20834812
$ws.Cells.Item(122, 9).Value = 1593.1
$ws.Cells.Item(122, 10).Value = 125000904
$ws.Cells.Item(122, 11).Value = 4779.299999999999
$ws.Cells.Item(122, 12).Value = 375002712
$ws.Cells.Item(122, 13).Value = -2329.299999999999
$ws.Cells.Item(122, 14).Value = -375007612

$ws.Cells.Item(132, 8).Value = 1998.1714
$ws.Cells.Item(132, 9).Value = 1764.2667
$ws.Cells.Item(132, 10).Value = 3401.6
$ws.Cells.Item(132, 11).Value = 5292.800099999999
$ws.Cells.Item(132, 12).Value = 10204.8
$ws.Cells.Item(132, 13).Value = -2762.800099999999
$ws.Cells.Item(132, 14).Value = -15264.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(29, 8).Value = 100000000
$ws.Cells.Item(29, 10).Value = 100000000
$ws.Cells.Item(29, 12).Value = 100000000
$ws.Cells.Item(29, 14).Value = -100000590

$ws.Cells.Item(34, 8).Value = 10000
$ws.Cells.Item(34, 9).Value = 10000
$ws.Cells.Item(34, 11).Value = 10000
$ws.Cells.Item(34, 13).Value = -9828

$ws.Cells.Item(40, 8).Value = 145372
$ws.Cells.Item(40, 9).Value = 169117.33
$ws.Cells.Item(40, 10).Value = 2900
$ws.Cells.Item(40, 11).Value = 169117.33
$ws.Cells.Item(40, 12).Value = 2900
$ws.Cells.Item(40, 13).Value = -168981.33
$ws.Cells.Item(40, 14).Value = -3172

$ws.Cells.Item(93, 8).Value = 6769.222
$ws.Cells.Item(93, 9).Value = 8396.714
$ws.Cells.Item(93, 10).Value = 1073
$ws.Cells.Item(93, 11).Value = 8396.714
$ws.Cells.Item(93, 12).Value = 1073
$ws.Cells.Item(93, 13).Value = -7148.714
$ws.Cells.Item(93, 14).Value = -3569

$ws.Cells.Item(132, 8).Value = 2240.5557
$ws.Cells.Item(132, 9).Value = 1509.4783
$ws.Cells.Item(132, 11).Value = 4528.4349
$ws.Cells.Item(132, 13).Value = -1998.4349

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(80, 8).Value = 62250
$ws.Cells.Item(80, 10).Value = 62250
$ws.Cells.Item(80, 12).Value = 62250
$ws.Cells.Item(80, 14).Value = -64246

$ws.Cells.Item(83, 8).Value = 62250
$ws.Cells.Item(83, 10).Value = 62250
$ws.Cells.Item(83, 12).Value = 186750
$ws.Cells.Item(83, 14).Value = -196734

$ws.Cells.Item(122, 8).Value = 2407.9
$ws.Cells.Item(122, 9).Value = 2347.375
$ws.Cells.Item(122, 10).Value = 2650
$ws.Cells.Item(122, 11).Value = 7042.125
$ws.Cells.Item(122, 12).Value = 7950
$ws.Cells.Item(122, 13).Value = -4592.125
$ws.Cells.Item(122, 14).Value = -12850

$ws.Cells.Item(132, 8).Value = 911.725
$ws.Cells.Item(132, 9).Value = 644.0645
$ws.Cells.Item(132, 11).Value = 1932.1935
$ws.Cells.Item(132, 13).Value = 597.8065000000001
